$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove C2, E2, C3 (now empty cells)
$ws.Range("C2").Value = $null
$ws.Range("E2").Value = $null
$ws.Range("C3").Value = $null

# Update remaining cells with corrected (refined) forecast values
$ws.Range("E3").Value = 1.087227286828263
$ws.Range("C4").Value = -4.774178217057756
$ws.Range("E4").Value = -0.8523446516643496
$ws.Range("E5").Value = 1.390521443873438
$ws.Range("C7").Value = 1.239479831392831
$ws.Range("C8").Value = 0.2379616621360992
$ws.Range("E9").Value = 0.6176326357195894
$ws.Range("C10").Value = 1.470039379455734
$ws.Range("E10").Value = 1.577608035818323
$ws.Range("C11").Value = 1.638797242243228
$ws.Range("E11").Value = 1.369334405341593
$ws.Range("E13").Value = 1.36203066512679
$ws.Range("E15").Value = -0.1895486537906388
$ws.Range("C16").Value = 1.099928004397577
$ws.Range("C17").Value = 2.310042359896247
$ws.Range("E18").Value = 1.063035646777677
$ws.Range("C19").Value = -0.3101476031197037
